# Apply updated dSF (column F) values per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -4
$ws.Range("F8").Value = -3
$ws.Range("F10").Value = 7
$ws.Range("F13").Value = 13
$ws.Range("F16").Value = -5
$ws.Range("F18").Value = 0
$ws.Range("F21").Value = -3
